$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column G: "wanderingMarker" trigger (trigger 13 for mindwandering) ---
$ws.Range("G1").Value = "wanderingMarker"
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 7).Value = 13
}

# Column widths for F and G (best-fit sized to their new/longer content)
$ws.Columns.Item(6).ColumnWidth = 11.6666666667
$ws.Columns.Item(7).ColumnWidth = 16

# --- View adjustments ---
$win = $excel.ActiveWindow
$win.Zoom = 100

# Re-seat the frozen top row and scroll the body down near the new data,
# then leave the new last cell selected.
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$win.ScrollRow = 23

$ws.Range("G42").Select()
